$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("A1").Value = 1111111111
$ws.Range("A2").Value = 222222222
$ws.Range("A3").Value = 333333
$ws.Range("A4").Value = 4444444
$ws.Range("A5").Value = 5555555

# Best-fit column A to its new (wider) content, matching
# <col min="1" max="1" width="11" bestFit="1" customWidth="1"/>
$ws.Columns("A").AutoFit()
$ws.Columns("A").ColumnWidth = 10.1667

# Update selection to F6
$ws.Range("F6").Select()
